$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update match-id (column A) and numeric matchweek (column E) for all data rows ---
$ws.Range("A2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("A3").Value = 4
$ws.Range("E3").Value = 3
$ws.Range("A4").Value = 7
$ws.Range("E4").Value = 5
$ws.Range("A5").Value = 9
$ws.Range("E5").Value = 7
$ws.Range("A6").Value = 11
$ws.Range("E6").Value = 8
$ws.Range("A7").Value = 15
$ws.Range("E7").Value = 11
$ws.Range("A8").Value = 16
$ws.Range("E8").Value = 12
$ws.Range("A9").Value = 19
$ws.Range("E9").Value = 14
$ws.Range("A10").Value = 21
$ws.Range("E10").Value = 16
$ws.Range("A11").Value = 26
$ws.Range("E11").Value = 19
$ws.Range("A12").Value = 28
$ws.Range("E12").Value = 21
$ws.Range("A13").Value = 30
$ws.Range("E13").Value = 23
$ws.Range("A14").Value = 32
$ws.Range("E14").Value = 25
$ws.Range("A15").Value = 21
$ws.Range("E15").Value = 15
$ws.Range("A16").Value = 3
$ws.Range("E16").Value = 4
$ws.Range("A17").Value = 8
$ws.Range("E17").Value = 6
$ws.Range("A18").Value = 26
$ws.Range("E18").Value = 17
$ws.Range("A19").Value = 12
$ws.Range("E19").Value = 9
$ws.Range("A20").Value = 33
$ws.Range("E20").Value = 26
$ws.Range("A21").Value = 18
$ws.Range("E21").Value = 10
$ws.Range("A22").Value = 1
$ws.Range("E22").Value = 2
$ws.Range("A25").Value = 19
$ws.Range("E25").Value = 20
$ws.Range("A26").Value = 17
$ws.Range("E26").Value = 18
$ws.Range("A27").Value = 34
$ws.Range("E27").Value = 24

# --- Rows 23 and 24: swap all match-stat columns (B:BD, excluding A/E handled below) ---
$ws.Range("B23").Value = "'2023-09-03"
$ws.Range("B24").Value = "'2023-07-02"
$ws.Range("C23").Value = "18:30"
$ws.Range("C24").Value = "11:00"
$ws.Range("D23").Value = "Série A"
$ws.Range("D24").Value = "Série A"
$ws.Range("F23").Value = "Sun"
$ws.Range("F24").Value = "Sun"
$ws.Range("G23").Value = "Home"
$ws.Range("G24").Value = "Home"
$ws.Range("H23").Value = "D"
$ws.Range("H24").Value = "L"
$ws.Range("I23").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("J24").Value = 1
$ws.Range("K23").Value = "Bragantino"
$ws.Range("K24").Value = "Bragantino"
$ws.Range("L23").Value = 1.2
$ws.Range("L24").Value = 0.6
$ws.Range("M23").Value = 1.2
$ws.Range("M24").Value = 0.9
$ws.Range("N23").Value = 50
$ws.Range("N24").Value = 42
$ws.Range("O23").ClearContents()
$ws.Range("O24").Value = 42404
$ws.Range("P23").Value = 15
$ws.Range("P24").Value = 12
$ws.Range("Q23").Value = 3
$ws.Range("Q24").Value = 4
$ws.Range("R23").Value = 20
$ws.Range("R24").Value = 33.3
$ws.Range("S23").Value = 0
$ws.Range("S24").Value = 0
$ws.Range("T23").Value = 0
$ws.Range("T24").Value = 0
$ws.Range("U23").Value = 0
$ws.Range("U24").Value = 1
$ws.Range("V23").Value = 0
$ws.Range("V24").Value = 0
$ws.Range("W23").Value = 0
$ws.Range("W24").Value = 0
$ws.Range("X23").Value = 1.2
$ws.Range("X24").Value = 0.6
$ws.Range("Y23").Value = 0.08
$ws.Range("Y24").Value = 0.05
$ws.Range("Z23").Value = -1.2
$ws.Range("Z24").Value = -0.6
$ws.Range("AA23").Value = -1.2
$ws.Range("AA24").Value = -0.6
$ws.Range("AB23").Value = 5
$ws.Range("AB24").Value = 5
$ws.Range("AC23").Value = 5
$ws.Range("AC24").Value = 4
$ws.Range("AD23").Value = 100
$ws.Range("AD24").Value = 80
$ws.Range("AE23").Value = 1
$ws.Range("AE24").Value = 0
$ws.Range("AF23").Value = 1.1
$ws.Range("AF24").Value = 1.4
$ws.Range("AG23").Value = 1.1
$ws.Range("AG24").Value = 0.4
$ws.Range("AH23").Value = 5049
$ws.Range("AH24").Value = 5119
$ws.Range("AI23").Value = 2385
$ws.Range("AI24").Value = 2197
$ws.Range("AJ23").Value = 0
$ws.Range("AJ24").Value = 0
$ws.Range("AK23").Value = 1.1
$ws.Range("AK24").Value = 0.4
$ws.Range("AL23").Value = 1.8
$ws.Range("AL24").Value = 0.5
$ws.Range("AM23").Value = 11
$ws.Range("AM24").Value = 9
$ws.Range("AN23").Value = 27
$ws.Range("AN24").Value = 16
$ws.Range("AO23").Value = 9
$ws.Range("AO24").Value = 5
$ws.Range("AP23").Value = 1
$ws.Range("AP24").Value = 1
$ws.Range("AQ23").Value = 31
$ws.Range("AQ24").Value = 24
$ws.Range("AR23").Value = 3
$ws.Range("AR24").Value = 2
$ws.Range("AS23").Value = 1
$ws.Range("AS24").Value = 6
$ws.Range("AT23").Value = 30
$ws.Range("AT24").Value = 17
$ws.Range("AU23").Value = 17
$ws.Range("AU24").Value = 14
$ws.Range("AV23").Value = 24
$ws.Range("AV24").Value = 23
$ws.Range("AW23").Value = 0
$ws.Range("AW24").Value = 0
$ws.Range("AX23").Value = 11
$ws.Range("AX24").Value = 13
$ws.Range("AY23").Value = 11
$ws.Range("AY24").Value = 9
$ws.Range("AZ23").Value = 7
$ws.Range("AZ24").Value = 11
$ws.Range("BA23").Value = 1
$ws.Range("BA24").Value = 0
$ws.Range("BB23").Value = 8
$ws.Range("BB24").Value = 4
$ws.Range("BC23").Value = 1
$ws.Range("BC24").Value = 0
$ws.Range("BD23").Value = "Cruzeiro"
$ws.Range("BD24").Value = "Corinthians"

# --- Rows 23/24 new match-id (A) and matchweek (E) ---
$ws.Range("A23").Value = 21
$ws.Range("E23").Value = 22
$ws.Range("A24").Value = 18
$ws.Range("E24").Value = 13
